# Add the new "2022-Q3" sheet, positioned right after "总计" and before "2022-Q2",
# populate it with fund-holding data, and update the "总计" summary sheet with a
# new leading row for 2022-Q3 (pushing the existing rows down).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new quarter worksheet --------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

# Header row (bold / centered, matches the other quarter sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data rows
$rows = @(
    @(0, "005585", "银河文体娱乐主题灵活配置混合A", "3.01", "90.28", "5.23", "0.1574", 4),
    @(1, "001628", "招商体育文化休闲股票A",          "2.23", "92.42", "5.10", "0.1137", 4),
    @(2, "015667", "银河文体娱乐主题灵活配置混合C", "0.41", "90.28", "5.23", "0.0214", 4),
    @(3, "015395", "招商体育文化休闲股票C",          "0.25", "92.42", "5.10", "0.0128", 4),
    @(4, "001735", "广发百发大数据策略成长灵活配置混合E", "0.51", "43.77", "0.84", "0.0043", 6),
    @(5, "001734", "广发百发大数据策略成长灵活配置混合A", "0.27", "43.77", "0.84", "0.0023", 6)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[1]

    $q3.Cells.Item($r, 3).Value = $row[2]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[3]

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[4]

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[5]

    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[6]

    $q3.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# --- 2. Update the "总计" summary sheet ----------------------------------
# Insert a new row above the current row 2 so the existing quarters shift
# down, then fill it in with the 2022-Q3 totals.
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.31

# Renumber the "序号"-style column A for the rows that got pushed down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
